# Reflect changes to in-progress User Stories on the Sprint Plan (Sheet1)
# and the Burndown Chart that is driven from it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "Planned" burndown values for 7/9 - 7/15 (rows 9-15) from 16 to 11.
$ws.Range("B9:B15").Value = 11

# Move the active selection to reflect where the user was last working.
$ws.Activate()
$ws.Range("C10").Select()
